# Clean up the "Reference" column (A) by stripping an erroneous trailing
# "16" that was appended to the end of each Bible reference (e.g.
# "2 Peter 1:116" -> "2 Peter 1:1"), to prepare the data for human
# readability. Rows that already look correct (no trailing "16") are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $value = $cell.Value2

    if ($value -ne $null -and $value.EndsWith("16")) {
        $cell.Value = $value.Substring(0, $value.Length - 2)
    }
}
